$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Linked missing student number data": the institutions list is missing the
# "University of X" variant name for three institutions that are already
# present under their short-form name. Add a row with the alternate name
# directly below each existing "<domain> | <Short Name> University" row,
# same as if the user selected that row and chose Insert > Copied/Entire Row.
# Work from the bottom of the sheet upwards so earlier (lower) row numbers
# used below are not shifted by a later insertion.

# University of St. Andrews -> right after row 69 (st-andrews | University of St Andrews)
$ws.Rows.Item(70).Insert()
$ws.Range("A70").Value = "st-andrews"
$ws.Range("B70").Value = "University of St. Andrews"

# University of Lancaster -> right after row 22 (lancaster | Lancaster University)
$ws.Rows.Item(23).Insert()
$ws.Range("A23").Value = "lancaster"
$ws.Range("B23").Value = "University of Lancaster"

# University of Durham -> right after row 15 (durham | Durham University)
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = "durham"
$ws.Range("B16").Value = "University of Durham"

# The table (and its AutoFilter / _FilterDatabase name) now spans 81 data
# rows instead of 78. Re-apply the same domain filter (unchanged set of
# values/unchecked items) over the new, larger range so the filter's
# remembered criteria - and the rows it hides - stay exactly as they were.
$ws.AutoFilterMode = $false
$domainFilterValues = @("aber","abertay","arts","aston","bath","bcu","bristol","brookes","cam","cardiff","chester","citystgeorges","coventry","dundee","durham","ed","falmouth","gcu","gla","glos","gold","herts","hud","hull","hw","kcl","kent","lancaster","lboro","leeds","leedstrinity","lincoln","liverpool","lse","manchester","mdx","mmu","napier","ncl","newman","northumbria","nottingham","open","ox","port","qmul","qub","rgu","roehampton","sheffield","shu","solent","southampton","southwales","st-andrews","stir","strath","sunderland","sussex","swansea","uclan","ulster","uwl","warwick","york","yorksj")
$ws.Range("A2:B81").AutoFilter(1, $domainFilterValues, 7)

# Keep the workbook-level hidden "_FilterDatabase" name in sync with the
# AutoFilter's new range, same as Excel does automatically.
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$2:`$B`$81"

# Reflect where editing ended up (last new row added).
$ws.Range("B72").Select()
